$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7724602222442627
$ws.Range("B1").Value = 1.030999064445496
$ws.Range("C1").Value = 3.692499876022339
$ws.Range("D1").Value = 2.180202484130859
$ws.Range("E1").Value = 1.093103408813477
